# Update the fixed "footer" date shown on the slide master and every slide
# layout from 2022. 4. 28. to 2022. 4. 29. (Insert > Header & Footer > Date
# fixed text, applied to all masters/layouts).

$p = $ppt.ActivePresentation

$oldDate = "2022. 4. 28."
$newDate = "2022. 4. 29."

# ppPlaceholderDate
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if (-not $shape.HasTextFrame) { continue }

        $isDatePlaceholder = $false
        if ($shape.Type -eq 14) {
            if ($shape.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        }

        if (-not $isDatePlaceholder) { continue }

        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout hanging off the master
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
